$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = "[4.3619097694546785, 8.805726194114495]"
$ws.Range("M2").Value = 0.00000001127610782525323
$ws.Range("N2").Value = 0.00000001127610782525323
$ws.Range("P2").Value = "[-1.7610529389492324, -1.0063159651138474]"
$ws.Range("Q2").Value = 0.000000000002597699833017941
$ws.Range("R2").Value = 0.000000000005195399666035883
$ws.Range("T2").Value = "[8.019973099640499, 10.596555111024475]"
$ws.Range("X2").Value = 3.666066066066087
$ws.Range("Y2").Value = 6.41561561561565

# Row 3
$ws.Range("L3").Value = "[5.082948930786207, 8.834268073268845]"
$ws.Range("M3").Value = 0.000000000001297628671181883
$ws.Range("N3").Value = 0.000000000002595257342363766
$ws.Range("P3").Value = "[-0.2327105669325764, 0.39623691126357663]"
$ws.Range("Q3").Value = 0.6096665997523454
$ws.Range("R3").Value = 0.6096665997523454
$ws.Range("T3").Value = "[7.79031577900189, 10.16604868233166]"
$ws.Range("X3").Value = 21.92432432432453
$ws.Range("Y3").Value = 24.26666666666689
